$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "51.593.89"
$ws.Range("E2").Value = "  +0.04%  "
Set-TextValue $ws.Range("D3") "2.790.92"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "352.94"
$ws.Range("E5").Value = "  -1.59%  "
Set-TextValue $ws.Range("D6") "111.11"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue $ws.Range("D9") "0.628"
$ws.Range("E9").Value = "  +7.05%  "
Set-TextValue $ws.Range("D10") "39.99"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("E11").Value = "  -1.97%  "
Set-TextValue $ws.Range("D12") "0.0836"
$ws.Range("E12").Value = "  -0.45%  "
Set-TextValue $ws.Range("D13") "19.91"
$ws.Range("E13").Value = "  +1.59%  "
Set-TextValue $ws.Range("D14") "7.73"
$ws.Range("E14").Value = "  +1.98%  "
Set-TextValue $ws.Range("D15") "3.233.41"
$ws.Range("E15").Value = "  +0.72%  "
Set-TextValue $ws.Range("D16") "2.795.60"
$ws.Range("E16").Value = "  +0.12%  "
Set-TextValue $ws.Range("D17") "0.941"
$ws.Range("E17").Value = "  +1.89%  "
Set-TextValue $ws.Range("D18") "51.570.89"
$ws.Range("E18").Value = "  +0.03%  "
Set-TextValue $ws.Range("D19") "7.57"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").Value = "  +3.78%  "
Set-TextValue $ws.Range("D21") "13.52"
$ws.Range("E22").Value = "  +0.62%  "
Set-TextValue $ws.Range("D23") "70.20"
$ws.Range("E23").Value = "  +0.54%  "
Set-TextValue $ws.Range("D24") "266.87"
$ws.Range("E24").Value = "  -0.28%  "
Set-TextValue $ws.Range("D25") "2.74"
$ws.Range("E25").Value = "  -0.71%  "
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  +0.02%  "
Set-TextValue $ws.Range("D27") "26.03"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -0.98%  "
Set-TextValue $ws.Range("D29") "38.90"
$ws.Range("E29").Value = "  +10.94%  "
Set-TextValue $ws.Range("D30") "10.32"
$ws.Range("E30").Value = "  +1.99%  "
Set-TextValue $ws.Range("D31") "2.26"
$ws.Range("E31").Value = "  -1.91%  "
Set-TextValue $ws.Range("D32") "52.55"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("E33").Value = "  +1.12%  "
Set-TextValue $ws.Range("D34") "0.0453"
$ws.Range("E34").Value = "  +3.43%  "
Set-TextValue $ws.Range("D35") "0.0896"
$ws.Range("E35").Value = "  +7.15%  "
Set-TextValue $ws.Range("D36") "5.57"
$ws.Range("E36").Value = "  +8.53%  "
Set-TextValue $ws.Range("D37") "1.00"
$ws.Range("E37").Value = "  +0.03%  "
Set-TextValue $ws.Range("D38") "18.80"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +2.74%  "
$ws.Range("E41").Value = "  +0.86%  "
Set-TextValue $ws.Range("D42") "2.49"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +1.86%  "
Set-TextValue $ws.Range("D44") "121.24"
$ws.Range("E44").Value = "  +0.69%  "
Set-TextValue $ws.Range("D45") "21.76"
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("E46").Value = "  +6.65%  "
Set-TextValue $ws.Range("D47") "3.39"
$ws.Range("E47").Value = "  +4.43%  "
Set-TextValue $ws.Range("D48") "2.103.86"
$ws.Range("E48").Value = "  +1.10%  "
Set-TextValue $ws.Range("D49") "0.950"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D50") "5.45"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D51") "1.36"
$ws.Range("E51").Value = "  +6.24%  "
